$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 383.14285
$ws.Range("I33").Value = 338.4
$ws.Range("J33").Value = 495
$ws.Range("K33").Value = 338.4
$ws.Range("L33").Value = 495
$ws.Range("M33").Value = -109.4
$ws.Range("N33").Value = -953

$ws.Range("H125").Value = 782.75
$ws.Range("I125").Value = 683.6667
$ws.Range("J125").Value = 881.8333
$ws.Range("K125").Value = 6153.0003
$ws.Range("L125").Value = 7936.4997
$ws.Range("M125").Value = -3693.0003
$ws.Range("N125").Value = -12856.4997

$ws.Range("H129").Value = 1952.0605
$ws.Range("J129").Value = 2154.3103
$ws.Range("L129").Value = 6462.9309
$ws.Range("N129").Value = -16462.9309

$ws.Range("H131").Value = 2074.7
$ws.Range("I131").Value = 792.4167
$ws.Range("J131").Value = 3998.125
$ws.Range("K131").Value = 2377.2501
$ws.Range("L131").Value = 11994.375
$ws.Range("M131").Value = 2662.7499
$ws.Range("N131").Value = -22074.375

$ws.Range("H132").Value = 33035090
$ws.Range("I132").Value = 42252564
$ws.Range("K132").Value = 126757692
$ws.Range("M132").Value = -126755162

$ws.Range("H137").Value = 713764.9399999999
$ws.Range("I137").Value = 1590407.8
$ws.Range("J137").Value = 2973.4055
$ws.Range("K137").Value = 4771223.4
$ws.Range("L137").Value = 8920.216499999999
$ws.Range("M137").Value = -4768673.4
$ws.Range("N137").Value = -14020.2165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4738.7964
$ws.Range("I32").Value = 4532.787
$ws.Range("K32").Value = 4532.787
$ws.Range("M32").Value = -4245.787

$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490

$ws.Range("H109").Value = 33200
$ws.Range("J109").Value = 33200
$ws.Range("L109").Value = 33200
$ws.Range("N109").Value = -35974

$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H119").Value = 35224
$ws.Range("J119").Value = 35224
$ws.Range("L119").Value = 35224
$ws.Range("N119").Value = -44900

$ws.Range("H121").Value = 27169
$ws.Range("J121").Value = 27169
$ws.Range("L121").Value = 27169
$ws.Range("N121").Value = -30663

$ws.Range("H137").Value = 39811
$ws.Range("J137").Value = 39811
$ws.Range("L137").Value = 39811
$ws.Range("N137").Value = -50011

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2408.2104
$ws.Range("I11").Value = 801.4
$ws.Range("J11").Value = 2982.0715
$ws.Range("K11").Value = 801.4
$ws.Range("L11").Value = 2982.0715
$ws.Range("M11").Value = -661.4
$ws.Range("N11").Value = -3262.0715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 551.24243
$ws.Range("I22").Value = 336.35
$ws.Range("J22").Value = 881.8461
$ws.Range("K22").Value = 336.35
$ws.Range("L22").Value = 881.8461
$ws.Range("M22").Value = 13.64999999999998
$ws.Range("N22").Value = -1581.8461

$ws.Range("H25").Value = 34998.75
$ws.Range("J25").Value = 34998.75
$ws.Range("L25").Value = 34998.75
$ws.Range("N25").Value = -35346.75

$ws.Range("H86").Value = 2129.96
$ws.Range("I86").Value = 2017
$ws.Range("J86").Value = 2330.7778
$ws.Range("K86").Value = 2017
$ws.Range("L86").Value = 2330.7778
$ws.Range("M86").Value = -894
$ws.Range("N86").Value = -4576.7778

$ws.Range("H89").Value = 2129.96
$ws.Range("I89").Value = 2017
$ws.Range("J89").Value = 2330.7778
$ws.Range("K89").Value = 10085
$ws.Range("L89").Value = 11653.889
$ws.Range("M89").Value = -4469
$ws.Range("N89").Value = -22885.889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2737.8484
$ws.Range("I68").Value = 1093.875
$ws.Range("J68").Value = 3263.92
$ws.Range("K68").Value = 3281.625
$ws.Range("L68").Value = 9791.76
$ws.Range("M68").Value = -2470.625
$ws.Range("N68").Value = -11413.76

$ws.Range("H71").Value = 2737.8484
$ws.Range("I71").Value = 1093.875
$ws.Range("J71").Value = 3263.92
$ws.Range("K71").Value = 9844.875
$ws.Range("L71").Value = 29375.28
$ws.Range("M71").Value = -5788.875
$ws.Range("N71").Value = -37487.28

$ws.Range("J76").Value = 3000
$ws.Range("L76").Value = 9000
$ws.Range("N76").Value = -9766

$ws.Range("J79").Value = 3000
$ws.Range("L79").Value = 9000
$ws.Range("N79").Value = -11652

$ws.Range("H86").Value = 954.8
$ws.Range("J86").Value = 1401.6666
$ws.Range("L86").Value = 4204.9998
$ws.Range("N86").Value = -6576.9998

$ws.Range("H89").Value = 954.8
$ws.Range("J89").Value = 1401.6666
$ws.Range("L89").Value = 12614.9994
$ws.Range("N89").Value = -24470.9994

$ws.Range("H113").Value = 1812288.4
$ws.Range("I113").Value = 633.1539
$ws.Range("J113").Value = 7353822
$ws.Range("K113").Value = 1899.4617
$ws.Range("L113").Value = 22061466
$ws.Range("M113").Value = 270.5382999999999
$ws.Range("N113").Value = -22065806

$ws.Range("H132").Value = 2171.9
$ws.Range("J132").Value = 3203.4546
$ws.Range("L132").Value = 28831.0914
$ws.Range("N132").Value = -33891.0914

$ws.Range("H137").Value = 3231.12
$ws.Range("I137").Value = 4017.2727
$ws.Range("J137").Value = 2613.4285
$ws.Range("K137").Value = 12051.8181
$ws.Range("L137").Value = 7840.2855
$ws.Range("M137").Value = -6951.8181
$ws.Range("N137").Value = -18040.2855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1251942.2

$ws.Range("H8").Value = 1251942.2

$ws.Range("H122").Value = 4430.4
$ws.Range("I122").Value = 2801
$ws.Range("J122").Value = 6874.5
$ws.Range("K122").Value = 8403
$ws.Range("L122").Value = 20623.5
$ws.Range("M122").Value = -5953
$ws.Range("N122").Value = -25523.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3528.4348
$ws.Range("I7").Value = 1474.9
$ws.Range("K7").Value = 1474.9
$ws.Range("M7").Value = -1362.9

$ws.Range("H40").Value = 5915
$ws.Range("I40").Value = 5255.4546
$ws.Range("K40").Value = 5255.4546
$ws.Range("M40").Value = -5119.4546

$ws.Range("H46").Value = 1005.6774
$ws.Range("I46").Value = 660.7619
$ws.Range("J46").Value = 1730
$ws.Range("K46").Value = 660.7619
$ws.Range("L46").Value = 1730
$ws.Range("M46").Value = -472.7619
$ws.Range("N46").Value = -2106

$ws.Range("H64").Value = 25000
$ws.Range("J64").Value = 25000
$ws.Range("L64").Value = 25000
$ws.Range("N64").Value = -25450

$ws.Range("H67").Value = 25000
$ws.Range("J67").Value = 25000
$ws.Range("L67").Value = 25000
$ws.Range("N67").Value = -26560

$ws.Range("H126").Value = 3528.4348
$ws.Range("I126").Value = 1474.9
$ws.Range("K126").Value = 4424.700000000001
$ws.Range("M126").Value = -1954.700000000001

$ws.Range("H132").Value = 3993.5938
$ws.Range("I132").Value = 2758.2354
$ws.Range("J132").Value = 5393.6665
$ws.Range("K132").Value = 8274.706200000001
$ws.Range("L132").Value = 16180.9995
$ws.Range("M132").Value = -5744.706200000001
$ws.Range("N132").Value = -21240.9995

$ws.Range("H136").Value = 2803.6086
$ws.Range("I136").Value = 1029.4615
$ws.Range("J136").Value = 5110
$ws.Range("K136").Value = 3088.3845
$ws.Range("L136").Value = 15330
$ws.Range("M136").Value = -538.3844999999997
$ws.Range("N136").Value = -20430

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 27849
$ws.Range("J119").Value = 27849
$ws.Range("L119").Value = 27849
$ws.Range("N119").Value = -37525

$ws.Range("H122").Value = 3191.75
$ws.Range("I122").Value = 1782.0588
$ws.Range("J122").Value = 4453.0527
$ws.Range("K122").Value = 5346.1764
$ws.Range("L122").Value = 13359.1581
$ws.Range("M122").Value = -2896.1764
$ws.Range("N122").Value = -18259.1581

$ws.Range("H126").Value = 446178.6
$ws.Range("I126").Value = 1545.091
$ws.Range("J126").Value = 822406.9399999999
$ws.Range("K126").Value = 4635.272999999999
$ws.Range("L126").Value = 2467220.82
$ws.Range("M126").Value = -2165.272999999999
$ws.Range("N126").Value = -2472160.82
